# Daily refresh of the cryptos list (GitHub Actions bot).
# Price (D) / Volume(1h) (E) columns are plain text in this sheet, and a
# couple of rows (48/49) swap rank between ARBITRUM and Stacks.
# Price values that look like a plain decimal number (e.g. "260.25") are
# prefixed with a leading apostrophe so Excel keeps storing them as text
# (matches how the sheet already stores values such as "44.195.02" that
# aren't valid numbers) instead of silently converting them to numerics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.241.70'
$ws.Range("E2").Value = '  +4.79%  '
$ws.Range("D3").Value = '2.227.67'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''260.25'
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").Value = '''83.44'
$ws.Range("E6").Value = '  +13.57%  '
$ws.Range("D7").Value = '''0.628'
$ws.Range("E7").Value = '  +4.08%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.610'
$ws.Range("E9").Value = '  +5.10%  '
$ws.Range("D10").Value = '''44.30'
$ws.Range("E10").Value = '  +10.88%  '
$ws.Range("E11").Value = '  +3.10%  '
$ws.Range("D12").Value = '''7.08'
$ws.Range("E12").Value = '  +4.83%  '
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("D14").Value = '2.560.30'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '''14.68'
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("D16").Value = '2.216.48'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("E17").Value = '  +2.83%  '
$ws.Range("D18").Value = '44.095.73'
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("D19").Value = '''0.0000105'
$ws.Range("E19").Value = '  +2.32%  '
$ws.Range("D20").Value = '''71.80'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("E21").Value = '  +3.48%  '
$ws.Range("D22").Value = '''2.36'
$ws.Range("E22").Value = '  +9.98%  '
$ws.Range("D23").Value = '''233.78'
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("E24").Value = '  -2.75%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '''10.80'
$ws.Range("E26").Value = '  +3.54%  '
$ws.Range("D27").Value = '''40.92'
$ws.Range("E27").Value = '  +11.28%  '
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '''2.25'
$ws.Range("E29").Value = '  +2.60%  '
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("D31").Value = '''173.16'
$ws.Range("E31").Value = '  +2.80%  '
$ws.Range("D32").Value = '''0.0894'
$ws.Range("E32").Value = '  +11.09%  '
$ws.Range("D33").Value = '''20.70'
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").Value = '''5.36'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  +9.27%  '
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("D37").Value = '''0.0368'
$ws.Range("E37").Value = '  +12.63%  '
$ws.Range("E38").Value = '  +7.31%  '
$ws.Range("D39").Value = '''13.40'
$ws.Range("E39").Value = '  +12.83%  '
$ws.Range("D40").Value = '''2.97'
$ws.Range("E40").Value = '  +23.98%  '
$ws.Range("E41").Value = '  +4.21%  '
$ws.Range("D42").Value = '''63.82'
$ws.Range("E42").Value = '  +8.60%  '
$ws.Range("D43").Value = '''5.56'
$ws.Range("E43").Value = '  +8.83%  '
$ws.Range("D44").Value = '''0.202'
$ws.Range("E44").Value = '  +3.32%  '
$ws.Range("D45").Value = '''103.23'
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D47").Value = '''8.38'
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''1.13'
$ws.Range("E48").Value = '  +4.20%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''1.57'
$ws.Range("E49").Value = '  +29.16%  '
$ws.Range("D50").Value = '''0.445'
$ws.Range("E50").Value = '  -3.40%  '
$ws.Range("E51").Value = '  +3.67%  '
